$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old placeholder content (Hej/sda/asd/sad) before laying out the new table.
$ws.Range("A1:C5").ClearContents()

# Header row (row 1) - League of Legends EU LCS game stats
$headers = @("League", "Date", "Blue Team", "Red Team", "FTBT", "FTRT", "FTR", "GBT", "GRT", "KBT", "KRT", "WBT", "WRT", "TDBT", "TDRT", "FBBT", "FBRT", "CWBT", "CWRT", "RBT", "RRT")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row (row 2)
$ws.Range("A2").Value = "EU LCS"
$ws.Range("B2").Value = 43176
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = "Splyce"
$ws.Range("D2").Value = "G2"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "R"
$ws.Range("H2").Value = 52000
$ws.Range("I2").Value = 64800
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 125
$ws.Range("M2").Value = 142
$ws.Range("N2").Value = 45400
$ws.Range("O2").Value = 64100
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 41
$ws.Range("S2").Value = 43

# Column widths approximating the bestFit autofit of the first four columns
$ws.Columns.Item(1).ColumnWidth = 6.5
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(4).ColumnWidth = 8.833333333333332

# Selection matching the target state
$ws.Range("U1").Select()
